# Apply the financial profile update across Summary, Assets and Liabilities
# sheets (borrower changed from Majid Al Memari to Noura Al Awani, along
# with refreshed balances pulled from the new data pipeline).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Noura Al Awani"
$summary.Range("B4").Value = 2509.36
$summary.Range("B6").Value = 857420
$summary.Range("B7").Value = 361519
$summary.Range("B8").Value = 495901
$summary.Range("B9").Value = 2.37

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")
$assets.Range("B2").Value = "Luxury Car"
$assets.Range("C2").Value = 575631
$assets.Range("B3").Value = "Premium Car"
$assets.Range("C3").Value = 280057
$assets.Range("C4").Value = 1732
$assets.Range("C5").Value = 857420

# ---------------------------------------------------------------------------
# Liabilities sheet
# ---------------------------------------------------------------------------
$liabilities = $wb.Worksheets.Item("Liabilities")

# Row 2 keeps its category/description (Auto Loans / Vehicle Loan 1) but the
# amount, monthly payment and remaining-years columns are refreshed.
$liabilities.Range("C2").Value = 345379
$liabilities.Range("D2").Value = 7195
$liabilities.Range("E2").Value = 4

# Row 3 becomes the Credit Cards / Credit Card Balance entry with new
# figures.
$liabilities.Range("A3").Value = "Credit Cards"
$liabilities.Range("B3").Value = "Credit Card Balance"
$liabilities.Range("C3").Value = 16140
$liabilities.Range("D3").Value = 807
$liabilities.Range("E3").Value = 1

# The old Personal Loans (row 4) and Credit Cards (row 5) rows are removed
# entirely, shifting the TOTAL LIABILITIES row up from row 6 to row 4.
$liabilities.Rows.Item(4).Delete()
$liabilities.Rows.Item(4).Delete()

# Refresh the TOTAL LIABILITIES figure (now sitting on row 4).
$liabilities.Range("C4").Value = 361519
